$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("I would like to have my hair cut.", "カットをお願いします。|カットをおねがいします。")
    ,@("Please don't make it too short.", "あまり短くしないでください。|あまりみじかくしないでください。")
    ,@("Please don't shave me.", "そらないでください。")
    ,@("Please cut off about 3 centimeters.", "３センチぐらい切ってください。|３センチぐらいきってください。")
    ,@("Please cut the back all the same length.", "後ろをそろえてください。|うしろをそろえてください。")
    ,@("Please dye my hair red.", "赤にそめてください。|あかにそめてください。")
    ,@("I want my hair to be like Bob Marley's.", "ボブ・マーリーみたいな髪形にしたいんですが。|ボブ・マーリーみたいなかみがたにしたいんですが。")
    ,@("shampoo", "シャンプー")
    ,@("cut", "カット")
    ,@("blow-dry", "ブロー")
    ,@("perm", "パーマ")
    ,@("hair coloring", "カラー")
    ,@("set", "セット")
    ,@("hair style", "髪形|かみがた")
    ,@("to cut", "切る|きる")
    ,@("to shave", "そる")
    ,@("to crop", "刈る|かる")
    ,@("to dye", "そめる")
    ,@("to make hair even; to trim", "そろえる")
    ,@("to have one's hair permed", "パーマをかける")
    ,@("to thin out (hair)", "すく")
    ,@("parting (of the hair)", "分け目|わけめ")
    ,@("bangs", "前髪|まえがみ")
    ,@("side", "横|よこ")
    ,@("back", "後ろ|うしろ")
)

$startRow = 60
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
